$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.328.33'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '3.043.03'
$ws.Range("E3").Value = '  +3.75%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '198.39'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.26'
$ws.Range("E6").Value = '  +3.99%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.208'
$ws.Range("E9").Value = '  +7.12%  '
$ws.Range("D10").Value = '3.041.41'
$ws.Range("E10").Value = '  +3.84%  '
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.19'
$ws.Range("E13").Value = '  +6.16%  '
$ws.Range("D14").Value = '3.599.25'
$ws.Range("E14").Value = '  +3.74%  '
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("D16").Value = '76.239.62'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").Value = '3.037.27'
$ws.Range("E18").Value = '  +3.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.52'
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.97'
$ws.Range("E20").Value = '  +4.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.64'
$ws.Range("E21").Value = '  +2.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.43'
$ws.Range("E22").Value = '  +7.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.39'
$ws.Range("E23").Value = '  +1.77%  '
$ws.Range("E24").Value = '  +3.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.44'
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.33'
$ws.Range("E27").Value = '  +3.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.84'
$ws.Range("E28").Value = '  +1.73%  '
$ws.Range("E29").Value = '  +1.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.995'
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.29'
$ws.Range("E31").Value = '  +2.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '492.36'
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.92'
$ws.Range("E34").Value = '  +5.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("E36").Value = '  +12.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.62'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.50'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.06'
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '192.60'
$ws.Range("E40").Value = '  +7.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.377'
$ws.Range("E41").Value = '  -4.07%  '
$ws.Range("E42").Value = '  -7.92%  '
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.790'
$ws.Range("E44").Value = '  +20.41%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.11'
$ws.Range("E45").Value = '  +4.96%  '
$ws.Range("E46").Value = '  +6.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '41.13'
$ws.Range("E47").Value = '  +2.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.64'
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.43'
$ws.Range("E49").Value = '  +6.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.594'
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.86'
$ws.Range("E51").Value = '  +0.00%  '
